$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2426969.2
$ws.Range("J17").Value = 2471866.8
$ws.Range("L17").Value = 7415600.399999999
$ws.Range("N17").Value = -7415936.399999999
$ws.Range("H18").Value = 385.5
$ws.Range("I18").Value = 385.5
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 385.5
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -101.5
$ws.Range("N18").Value = ""
$ws.Range("H33").Value = 1946.3182
$ws.Range("I33").Value = 262
$ws.Range("K33").Value = 262
$ws.Range("M33").Value = -33
$ws.Range("H40").Value = 2590.1904
$ws.Range("I40").Value = 2780.5625
$ws.Range("J40").Value = 1981
$ws.Range("K40").Value = 2780.5625
$ws.Range("L40").Value = 1981
$ws.Range("M40").Value = -2605.5625
$ws.Range("N40").Value = -2331
$ws.Range("H43").Value = 7165.8335
$ws.Range("J43").Value = 8749.25
$ws.Range("L43").Value = 8749.25
$ws.Range("N43").Value = -8887.25
$ws.Range("H74").Value = 3992.25
$ws.Range("I74").Value = 3994.5
$ws.Range("K74").Value = 3994.5
$ws.Range("M74").Value = -3058.5
$ws.Range("H76").Value = 5079.8335
$ws.Range("I76").Value = 4884.8335
$ws.Range("J76").Value = 5274.8335
$ws.Range("K76").Value = 4884.8335
$ws.Range("L76").Value = 5274.8335
$ws.Range("M76").Value = -4569.8335
$ws.Range("N76").Value = -5904.8335
$ws.Range("H77").Value = 3992.25
$ws.Range("I77").Value = 3994.5
$ws.Range("K77").Value = 19972.5
$ws.Range("M77").Value = -15292.5
$ws.Range("H79").Value = 5079.8335
$ws.Range("I79").Value = 4884.8335
$ws.Range("J79").Value = 5274.8335
$ws.Range("K79").Value = 4884.8335
$ws.Range("L79").Value = 5274.8335
$ws.Range("M79").Value = -3792.8335
$ws.Range("N79").Value = -7458.8335
$ws.Range("H92").Value = 3964.5715
$ws.Range("I92").Value = 4100.636
$ws.Range("J92").Value = 3465.6667
$ws.Range("K92").Value = 4100.636
$ws.Range("L92").Value = 3465.6667
$ws.Range("M92").Value = -2852.636
$ws.Range("N92").Value = -5961.6667
$ws.Range("H97").Value = 1430498.2
$ws.Range("J97").Value = 1430498.2
$ws.Range("L97").Value = 4291494.6
$ws.Range("N97").Value = -4292486.6
$ws.Range("H99").Value = 375362.25
$ws.Range("I99").Value = 167086.17
$ws.Range("K99").Value = 501258.51
$ws.Range("M99").Value = -499760.51
$ws.Range("H106").Value = 7996.625
$ws.Range("I106").Value = 1996.2858
$ws.Range("K106").Value = 1996.2858
$ws.Range("M106").Value = -1365.2858
$ws.Range("H113").Value = 12168.333
$ws.Range("I113").Value = 16501
$ws.Range("J113").Value = 3503
$ws.Range("K113").Value = 16501
$ws.Range("L113").Value = 3503
$ws.Range("M113").Value = -13247
$ws.Range("N113").Value = -10011
$ws.Range("H121").Value = 5795.6665
$ws.Range("J121").Value = 5795.6665
$ws.Range("L121").Value = 17386.9995
$ws.Range("N121").Value = -20880.9995
$ws.Range("H125").Value = 1940
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1940
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 17460
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -22380
$ws.Range("H131").Value = 14961.223
$ws.Range("I131").Value = 14518.875
$ws.Range("J131").Value = 18500
$ws.Range("K131").Value = 43556.625
$ws.Range("L131").Value = 55500
$ws.Range("M131").Value = -38516.625
$ws.Range("N131").Value = -65580
$ws.Range("H132").Value = 4580.6523
$ws.Range("I132").Value = 2991.0833
$ws.Range("K132").Value = 8973.249899999999
$ws.Range("M132").Value = -6443.249899999999
$ws.Range("H135").Value = 1855.4073
$ws.Range("I135").Value = 1533.5714
$ws.Range("K135").Value = 13802.1426
$ws.Range("M135").Value = -11267.1426
$ws.Range("H137").Value = 24433.963
$ws.Range("I137").Value = 30060.154
$ws.Range("J137").Value = 8761
$ws.Range("K137").Value = 90180.462
$ws.Range("L137").Value = 26283
$ws.Range("M137").Value = -87630.462
$ws.Range("N137").Value = -31383
$ws.Range("H138").Value = 2748.9883
$ws.Range("I138").Value = 1359.75
$ws.Range("K138").Value = 4079.25
$ws.Range("M138").Value = 1060.75
$ws.Range("H141").Value = 1143.4615
$ws.Range("I141").Value = 1143.4615
$ws.Range("K141").Value = 3430.3845
$ws.Range("M141").Value = 1749.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2066.6667
$ws.Range("I2").Value = 2066.6667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2066.6667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1953.6667
$ws.Range("N2").Value = ""
$ws.Range("H32").Value = 6190.439
$ws.Range("I32").Value = 3288.2122
$ws.Range("K32").Value = 3288.2122
$ws.Range("M32").Value = -3001.2122
$ws.Range("H61").Value = 3439.8293
$ws.Range("I61").Value = 2121.923
$ws.Range("J61").Value = 5724.2
$ws.Range("K61").Value = 2121.923
$ws.Range("L61").Value = 5724.2
$ws.Range("M61").Value = -1909.923
$ws.Range("N61").Value = -6148.2
$ws.Range("H74").Value = 63526.5
$ws.Range("I74").Value = 83800
$ws.Range("J74").Value = 2706
$ws.Range("K74").Value = 83800
$ws.Range("L74").Value = 2706
$ws.Range("M74").Value = -82926
$ws.Range("N74").Value = -4454
$ws.Range("H77").Value = 63526.5
$ws.Range("I77").Value = 83800
$ws.Range("J77").Value = 2706
$ws.Range("K77").Value = 419000
$ws.Range("L77").Value = 13530
$ws.Range("M77").Value = -414632
$ws.Range("N77").Value = -22266
$ws.Range("H97").Value = 930.2381
$ws.Range("I97").Value = 721.6667
$ws.Range("K97").Value = 721.6667
$ws.Range("M97").Value = -225.6667
$ws.Range("H116").Value = 2066.6667
$ws.Range("I116").Value = 2066.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2066.6667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 227.3332999999998
$ws.Range("N116").Value = ""
$ws.Range("H122").Value = 3051.7368
$ws.Range("I122").Value = 2374.5
$ws.Range("K122").Value = 7123.5
$ws.Range("M122").Value = -4673.5
$ws.Range("H132").Value = 1843.5588
$ws.Range("I132").Value = 1809.52
$ws.Range("J132").Value = 1938.1111
$ws.Range("K132").Value = 5428.559999999999
$ws.Range("L132").Value = 5814.3333
$ws.Range("M132").Value = -2898.559999999999
$ws.Range("N132").Value = -10874.3333
$ws.Range("H136").Value = 3439.8293
$ws.Range("I136").Value = 2121.923
$ws.Range("J136").Value = 5724.2
$ws.Range("K136").Value = 6365.768999999999
$ws.Range("L136").Value = 17172.6
$ws.Range("M136").Value = -3815.768999999999
$ws.Range("N136").Value = -22272.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2066.6667
$ws.Range("I3").Value = 2066.6667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2066.6667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1952.6667
$ws.Range("N3").Value = ""
$ws.Range("H20").Value = 2088.8696
$ws.Range("I20").Value = 2124.5
$ws.Range("J20").Value = 2033.4445
$ws.Range("K20").Value = 2124.5
$ws.Range("L20").Value = 2033.4445
$ws.Range("M20").Value = -1877.5
$ws.Range("N20").Value = -2527.4445
$ws.Range("H22").Value = 1428.875
$ws.Range("J22").Value = 2144
$ws.Range("L22").Value = 2144
$ws.Range("N22").Value = -2490
$ws.Range("H86").Value = 1543.8889
$ws.Range("I86").Value = 816.5
$ws.Range("J86").Value = 2998.6667
$ws.Range("K86").Value = 816.5
$ws.Range("L86").Value = 2998.6667
$ws.Range("M86").Value = 306.5
$ws.Range("N86").Value = -5244.6667
$ws.Range("H89").Value = 1543.8889
$ws.Range("I89").Value = 816.5
$ws.Range("J89").Value = 2998.6667
$ws.Range("K89").Value = 4082.5
$ws.Range("L89").Value = 14993.3335
$ws.Range("M89").Value = 1533.5
$ws.Range("N89").Value = -26225.3335
$ws.Range("H105").Value = 2101.6365
$ws.Range("I105").Value = 2051.8333
$ws.Range("K105").Value = 2051.8333
$ws.Range("M105").Value = -304.8332999999998
$ws.Range("H107").Value = 615.58826
$ws.Range("I107").Value = 495.60526
$ws.Range("K107").Value = 495.60526
$ws.Range("M107").Value = 1424.39474
$ws.Range("H132").Value = 103998.664
$ws.Range("J132").Value = 103998.664
$ws.Range("L132").Value = 103998.664
$ws.Range("N132").Value = -114118.664
$ws.Range("H134").Value = 3656.2812
$ws.Range("I134").Value = 2768.1
$ws.Range("J134").Value = 5136.5835
$ws.Range("K134").Value = 8304.299999999999
$ws.Range("L134").Value = 15409.7505
$ws.Range("M134").Value = -5769.299999999999
$ws.Range("N134").Value = -20479.7505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 111233224
$ws.Range("I4").Value = 13500
$ws.Range("K4").Value = 13500
$ws.Range("M4").Value = -13388
$ws.Range("H16").Value = 1586.4242
$ws.Range("I16").Value = 1504.1111
$ws.Range("J16").Value = 1685.2
$ws.Range("K16").Value = 1504.1111
$ws.Range("L16").Value = 1685.2
$ws.Range("M16").Value = -1217.1111
$ws.Range("N16").Value = -2259.2
$ws.Range("H22").Value = 367.5625
$ws.Range("I22").Value = 348.7143
$ws.Range("J22").Value = 499.5
$ws.Range("K22").Value = 348.7143
$ws.Range("L22").Value = 499.5
$ws.Range("M22").Value = 1.28570000000002
$ws.Range("N22").Value = -1199.5
$ws.Range("H31").Value = 335722.12
$ws.Range("I31").Value = 770383.1
$ws.Range("J31").Value = 3334.2942
$ws.Range("K31").Value = 770383.1
$ws.Range("L31").Value = 3334.2942
$ws.Range("M31").Value = -770088.1
$ws.Range("N31").Value = -3924.2942
$ws.Range("H34").Value = 335722.12
$ws.Range("I34").Value = 770383.1
$ws.Range("J34").Value = 3334.2942
$ws.Range("K34").Value = 770383.1
$ws.Range("L34").Value = 3334.2942
$ws.Range("M34").Value = -770181.1
$ws.Range("N34").Value = -3738.2942
$ws.Range("H58").Value = 3137.1667
$ws.Range("J58").Value = 2404.3333
$ws.Range("L58").Value = 2404.3333
$ws.Range("N58").Value = -2810.3333
$ws.Range("H60").Value = 11499
$ws.Range("I60").Value = 11499
$ws.Range("K60").Value = 11499
$ws.Range("M60").Value = -10988
$ws.Range("H62").Value = 3908.7
$ws.Range("I62").Value = 4081.6667
$ws.Range("J62").Value = 3649.25
$ws.Range("K62").Value = 4081.6667
$ws.Range("L62").Value = 3649.25
$ws.Range("M62").Value = -3457.6667
$ws.Range("N62").Value = -4897.25
$ws.Range("H65").Value = 3908.7
$ws.Range("I65").Value = 4081.6667
$ws.Range("J65").Value = 3649.25
$ws.Range("K65").Value = 20408.3335
$ws.Range("L65").Value = 18246.25
$ws.Range("M65").Value = -17288.3335
$ws.Range("N65").Value = -24486.25
$ws.Range("H68").Value = 39233.453
$ws.Range("J68").Value = 39233.453
$ws.Range("L68").Value = 39233.453
$ws.Range("N68").Value = -40731.453
$ws.Range("H71").Value = 39233.453
$ws.Range("J71").Value = 39233.453
$ws.Range("L71").Value = 117700.359
$ws.Range("N71").Value = -125188.359
$ws.Range("H94").Value = 1173.1072
$ws.Range("I94").Value = 693.36365
$ws.Range("J94").Value = 1483.5294
$ws.Range("K94").Value = 693.36365
$ws.Range("L94").Value = 1483.5294
$ws.Range("M94").Value = -242.36365
$ws.Range("N94").Value = -2385.5294
$ws.Range("H99").Value = 331046.22
$ws.Range("I99").Value = 724714.1
$ws.Range("J99").Value = 24860.055
$ws.Range("K99").Value = 724714.1
$ws.Range("L99").Value = 24860.055
$ws.Range("M99").Value = -723216.1
$ws.Range("N99").Value = -27856.055
$ws.Range("H105").Value = 5614.909
$ws.Range("I105").Value = 1776.25
$ws.Range("K105").Value = 1776.25
$ws.Range("M105").Value = -29.25
$ws.Range("H107").Value = 5859.069
$ws.Range("I107").Value = 1202.5714
$ws.Range("K107").Value = 1202.5714
$ws.Range("M107").Value = 717.4286
$ws.Range("H110").Value = 49999
$ws.Range("J110").Value = 49999
$ws.Range("L110").Value = 49999
$ws.Range("N110").Value = -58179
$ws.Range("H113").Value = 1586.4242
$ws.Range("I113").Value = 1504.1111
$ws.Range("J113").Value = 1685.2
$ws.Range("K113").Value = 1504.1111
$ws.Range("L113").Value = 1685.2
$ws.Range("M113").Value = 665.8888999999999
$ws.Range("N113").Value = -6025.2
$ws.Range("H122").Value = 3424.8572
$ws.Range("I122").Value = 3092
$ws.Range("J122").Value = 4257
$ws.Range("K122").Value = 9276
$ws.Range("L122").Value = 12771
$ws.Range("M122").Value = -6826
$ws.Range("N122").Value = -17671
$ws.Range("H126").Value = 331046.22
$ws.Range("I126").Value = 724714.1
$ws.Range("J126").Value = 24860.055
$ws.Range("K126").Value = 2174142.3
$ws.Range("L126").Value = 74580.16500000001
$ws.Range("M126").Value = -2171672.3
$ws.Range("N126").Value = -79520.16500000001
$ws.Range("H134").Value = 8456.762000000001
$ws.Range("I134").Value = 9033.166999999999
$ws.Range("J134").Value = 4998.3335
$ws.Range("K134").Value = 27099.501
$ws.Range("L134").Value = 14995.0005
$ws.Range("M134").Value = -24564.501
$ws.Range("N134").Value = -20065.0005
$ws.Range("H136").Value = 3137.1667
$ws.Range("J136").Value = 2404.3333
$ws.Range("L136").Value = 7212.999899999999
$ws.Range("N136").Value = -12312.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 2757.6
$ws.Range("I6").Value = 3862.6667
$ws.Range("J6").Value = 1100
$ws.Range("K6").Value = 11588.0001
$ws.Range("L6").Value = 3300
$ws.Range("M6").Value = -11475.0001
$ws.Range("N6").Value = -3526
$ws.Range("H7").Value = 50488.05
$ws.Range("I7").Value = 100306.9
$ws.Range("J7").Value = 669.2
$ws.Range("K7").Value = 300920.7
$ws.Range("L7").Value = 2007.6
$ws.Range("M7").Value = -300808.7
$ws.Range("N7").Value = -2231.6
$ws.Range("H17").Value = 7999
$ws.Range("J17").Value = 7999
$ws.Range("L17").Value = 23997
$ws.Range("N17").Value = -24335
$ws.Range("H37").Value = 77021090
$ws.Range("J37").Value = 77021090
$ws.Range("L37").Value = 231063270
$ws.Range("N37").Value = -231063494
$ws.Range("H38").Value = 59.857143
$ws.Range("J38").Value = 68.25
$ws.Range("L38").Value = 204.75
$ws.Range("N38").Value = -898.75
$ws.Range("H59").Value = 2500
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2500
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 7500
$ws.Range("M59").Value = ""
$ws.Range("N59").Value = -8580
$ws.Range("H97").Value = 235.16667
$ws.Range("J97").Value = 236
$ws.Range("L97").Value = 708
$ws.Range("N97").Value = -1700
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").Value = ""
$ws.Range("H121").Value = 130447.875
$ws.Range("I121").Value = 6997
$ws.Range("J121").Value = 148083.72
$ws.Range("K121").Value = 20991
$ws.Range("L121").Value = 444251.16
$ws.Range("M121").Value = -19681
$ws.Range("N121").Value = -446871.16
$ws.Range("H122").Value = 169
$ws.Range("I122").Value = 163.33333
$ws.Range("J122").Value = 170.54546
$ws.Range("K122").Value = 1469.99997
$ws.Range("L122").Value = 1534.90914
$ws.Range("M122").Value = 980.0000300000002
$ws.Range("N122").Value = -6434.90914
$ws.Range("H131").Value = 1390.3383
$ws.Range("I131").Value = 886.5
$ws.Range("K131").Value = 2659.5
$ws.Range("M131").Value = 2380.5
$ws.Range("H132").Value = 6985.909
$ws.Range("I132").Value = 8858.4375
$ws.Range("J132").Value = 1992.5
$ws.Range("K132").Value = 79725.9375
$ws.Range("L132").Value = 17932.5
$ws.Range("M132").Value = -77195.9375
$ws.Range("N132").Value = -22992.5
$ws.Range("H134").Value = 15039.417
$ws.Range("I134").Value = 15906.637
$ws.Range("K134").Value = 47719.911
$ws.Range("M134").Value = -42649.911

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 21960
$ws.Range("J43").Value = 29241.25
$ws.Range("L43").Value = 29241.25
$ws.Range("N43").Value = -29543.25
$ws.Range("H46").Value = 35409.832
$ws.Range("J46").Value = 45990
$ws.Range("L46").Value = 45990
$ws.Range("N46").Value = -46302
$ws.Range("H80").Value = 4640.619
$ws.Range("I80").Value = 3570.2856
$ws.Range("J80").Value = 5175.7856
$ws.Range("K80").Value = 3570.2856
$ws.Range("L80").Value = 5175.7856
$ws.Range("M80").Value = -2572.2856
$ws.Range("N80").Value = -7171.7856
$ws.Range("H83").Value = 4640.619
$ws.Range("I83").Value = 3570.2856
$ws.Range("J83").Value = 5175.7856
$ws.Range("K83").Value = 17851.428
$ws.Range("L83").Value = 25878.928
$ws.Range("M83").Value = -12859.428
$ws.Range("N83").Value = -35862.928
$ws.Range("H107").Value = 878.6923
$ws.Range("I107").Value = 1010.3333
$ws.Range("K107").Value = 1010.3333
$ws.Range("M107").Value = 909.6667
$ws.Range("H122").Value = 4900
$ws.Range("I122").Value = 4808.5264
$ws.Range("J122").Value = 5479.3335
$ws.Range("K122").Value = 14425.5792
$ws.Range("L122").Value = 16438.0005
$ws.Range("M122").Value = -11975.5792
$ws.Range("N122").Value = -21338.0005
$ws.Range("H132").Value = 36118.656
$ws.Range("I132").Value = 46254.043
$ws.Range("J132").Value = 5712.5
$ws.Range("K132").Value = 138762.129
$ws.Range("L132").Value = 17137.5
$ws.Range("M132").Value = -136232.129
$ws.Range("N132").Value = -22197.5
$ws.Range("H133").Value = 65999.57000000001
$ws.Range("J133").Value = 65999.57000000001
$ws.Range("L133").Value = 65999.57000000001
$ws.Range("N133").Value = -76119.57000000001
$ws.Range("H135").Value = 86999.336
$ws.Range("J135").Value = 86999.336
$ws.Range("L135").Value = 86999.336
$ws.Range("N135").Value = -97139.336
$ws.Range("H136").Value = 23517.334
$ws.Range("J136").Value = 23517.334
$ws.Range("L136").Value = 70552.00199999999
$ws.Range("N136").Value = -75652.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5744.077
$ws.Range("I7").Value = 5865.4165
$ws.Range("J7").Value = 4288
$ws.Range("K7").Value = 5865.4165
$ws.Range("L7").Value = 4288
$ws.Range("M7").Value = -5753.4165
$ws.Range("N7").Value = -4512
$ws.Range("H16").Value = 926.8
$ws.Range("I16").Value = 926.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 926.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -756.8
$ws.Range("N16").Value = ""
$ws.Range("H22").Value = 1188.6296
$ws.Range("I22").Value = 919.375
$ws.Range("J22").Value = 1302
$ws.Range("K22").Value = 919.375
$ws.Range("L22").Value = 1302
$ws.Range("M22").Value = -624.375
$ws.Range("N22").Value = -1892
$ws.Range("H27").Value = 1188.6296
$ws.Range("I27").Value = 919.375
$ws.Range("J27").Value = 1302
$ws.Range("K27").Value = 919.375
$ws.Range("L27").Value = 1302
$ws.Range("M27").Value = -812.375
$ws.Range("N27").Value = -1516
$ws.Range("H29").Value = 29999
$ws.Range("I29").Value = 29999
$ws.Range("K29").Value = 29999
$ws.Range("M29").Value = -29704
$ws.Range("H40").Value = 21791.23
$ws.Range("I40").Value = 26328.8
$ws.Range("K40").Value = 26328.8
$ws.Range("M40").Value = -26192.8
$ws.Range("H46").Value = 3208.6843
$ws.Range("I46").Value = 2988.3635
$ws.Range("J46").Value = 3511.625
$ws.Range("K46").Value = 2988.3635
$ws.Range("L46").Value = 3511.625
$ws.Range("M46").Value = -2800.3635
$ws.Range("N46").Value = -3887.625
$ws.Range("H55").Value = 5066.6855
$ws.Range("I55").Value = 767.1739
$ws.Range("J55").Value = 13307.417
$ws.Range("K55").Value = 767.1739
$ws.Range("L55").Value = 13307.417
$ws.Range("M55").Value = -594.1739
$ws.Range("N55").Value = -13653.417
$ws.Range("H61").Value = 1327.6522
$ws.Range("I61").Value = 1268.5294
$ws.Range("J61").Value = 1495.1666
$ws.Range("K61").Value = 1268.5294
$ws.Range("L61").Value = 1495.1666
$ws.Range("M61").Value = -1066.5294
$ws.Range("N61").Value = -1899.1666
$ws.Range("H68").Value = 3612.125
$ws.Range("I68").Value = 4249.5
$ws.Range("K68").Value = 4249.5
$ws.Range("M68").Value = -3500.5
$ws.Range("H71").Value = 3612.125
$ws.Range("I71").Value = 4249.5
$ws.Range("K71").Value = 21247.5
$ws.Range("M71").Value = -17503.5
$ws.Range("H113").Value = 1327.6522
$ws.Range("I113").Value = 1268.5294
$ws.Range("J113").Value = 1495.1666
$ws.Range("K113").Value = 1268.5294
$ws.Range("L113").Value = 1495.1666
$ws.Range("M113").Value = 901.4706000000001
$ws.Range("N113").Value = -5835.1666
$ws.Range("H122").Value = 4169.8
$ws.Range("I122").Value = 3667.3333
$ws.Range("J122").Value = 4923.5
$ws.Range("K122").Value = 11001.9999
$ws.Range("L122").Value = 14770.5
$ws.Range("M122").Value = -8551.999899999999
$ws.Range("N122").Value = -19670.5
$ws.Range("H126").Value = 5744.077
$ws.Range("I126").Value = 5865.4165
$ws.Range("J126").Value = 4288
$ws.Range("K126").Value = 17596.2495
$ws.Range("L126").Value = 12864
$ws.Range("M126").Value = -15126.2495
$ws.Range("N126").Value = -17804
$ws.Range("H132").Value = 2693.7666
$ws.Range("I132").Value = 2914.375
$ws.Range("J132").Value = 2613.5454
$ws.Range("K132").Value = 8743.125
$ws.Range("L132").Value = 7840.6362
$ws.Range("M132").Value = -6213.125
$ws.Range("N132").Value = -12900.6362
$ws.Range("H136").Value = 1114
$ws.Range("I136").Value = 1114
$ws.Range("K136").Value = 3342
$ws.Range("M136").Value = -792

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 19982.666
$ws.Range("J31").Value = 19982.666
$ws.Range("L31").Value = 19982.666
$ws.Range("N31").Value = -20678.666
$ws.Range("H51").Value = 20599.8
$ws.Range("H52").Value = 15599.4
$ws.Range("I52").Value = 13749.5
$ws.Range("K52").Value = 13749.5
$ws.Range("M52").Value = -13523.5
$ws.Range("H54").Value = 15722.5
$ws.Range("J54").Value = 15722.5
$ws.Range("L54").Value = 15722.5
$ws.Range("N54").Value = -16762.5
$ws.Range("H62").Value = 9031.454
$ws.Range("I62").Value = 7278.2
$ws.Range("K62").Value = 7278.2
$ws.Range("M62").Value = -6654.2
$ws.Range("H65").Value = 9031.454
$ws.Range("I65").Value = 7278.2
$ws.Range("K65").Value = 36391
$ws.Range("M65").Value = -33271
$ws.Range("H74").Value = 32198.8
$ws.Range("I74").Value = 17000
$ws.Range("J74").Value = 35998.5
$ws.Range("K74").Value = 17000
$ws.Range("L74").Value = 35998.5
$ws.Range("M74").Value = -16064
$ws.Range("N74").Value = -37870.5
$ws.Range("H77").Value = 32198.8
$ws.Range("I77").Value = 17000
$ws.Range("J77").Value = 35998.5
$ws.Range("K77").Value = 51000
$ws.Range("L77").Value = 107995.5
$ws.Range("M77").Value = -46320
$ws.Range("N77").Value = -117355.5
$ws.Range("H81").Value = 4238.5
$ws.Range("I81").Value = 4238.5
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8477
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -7416
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 4238.5
$ws.Range("I84").Value = 4238.5
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 42385
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -37081
$ws.Range("N84").Value = ""
$ws.Range("H109").Value = 24999.5
$ws.Range("J109").Value = 24999.5
$ws.Range("L109").Value = 24999.5
$ws.Range("N109").Value = -27773.5
$ws.Range("H113").Value = 957.9545000000001
$ws.Range("I113").Value = 811.5
$ws.Range("J113").Value = 1214.25
$ws.Range("K113").Value = 2434.5
$ws.Range("L113").Value = 3642.75
$ws.Range("M113").Value = -264.5
$ws.Range("N113").Value = -7982.75
$ws.Range("H122").Value = 20475.268
$ws.Range("I122").Value = 2086.75
$ws.Range("K122").Value = 6260.25
$ws.Range("M122").Value = -3810.25
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H126").Value = 3885
$ws.Range("I126").Value = 3885
$ws.Range("K126").Value = 11655
$ws.Range("M126").Value = -9185
$ws.Range("H132").Value = 2412.889
$ws.Range("I132").Value = 2391.4707
$ws.Range("K132").Value = 7174.4121
$ws.Range("M132").Value = -4644.4121
